# Auto-generated edit script applying market-data/profit updates
# across multiple sheets as described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2183.1667
$ws.Range("I29").Value = 100
$ws.Range("K29").Value = 300
$ws.Range("M29").Value = -19

$ws.Range("H43").Value = 5735.625
$ws.Range("J43").Value = 5314.1665
$ws.Range("L43").Value = 5314.1665
$ws.Range("N43").Value = -5452.1665

$ws.Range("H54").Value = 25000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 25000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 25000
$ws.Range("M54").Value = $null
$ws.Range("N54").Value = -25972

$ws.Range("H106").Value = 2615.8572
$ws.Range("I106").Value = 1904.3334
$ws.Range("K106").Value = 1904.3334
$ws.Range("M106").Value = -1273.3334

$ws.Range("H141").Value = 3613.3076
$ws.Range("I141").Value = 3470.524
$ws.Range("J141").Value = 4213
$ws.Range("K141").Value = 10411.572
$ws.Range("L141").Value = 12639
$ws.Range("M141").Value = -5231.572
$ws.Range("N141").Value = -22999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 9119.5
$ws.Range("I31").Value = 9355
$ws.Range("J31").Value = 7000
$ws.Range("K31").Value = 9355
$ws.Range("L31").Value = 7000
$ws.Range("M31").Value = -9061
$ws.Range("N31").Value = -7588

$ws.Range("H32").Value = 2200649.2
$ws.Range("I32").Value = 2697221.8
$ws.Range("K32").Value = 2697221.8
$ws.Range("M32").Value = -2696934.8

$ws.Range("H97").Value = 1635.875
$ws.Range("I97").Value = 1635.875
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1635.875
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1139.875
$ws.Range("N97").Value = $null

$ws.Range("H132").Value = 3579154.5
$ws.Range("I132").Value = 1916.375
$ws.Range("J132").Value = 10183287
$ws.Range("K132").Value = 5749.125
$ws.Range("L132").Value = 30549861
$ws.Range("M132").Value = -3219.125
$ws.Range("N132").Value = -30554921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null

$ws.Range("H99").Value = 18369.52
$ws.Range("J99").Value = 17129.572
$ws.Range("L99").Value = 17129.572
$ws.Range("N99").Value = -20125.572

$ws.Range("H134").Value = 36620.46
$ws.Range("I134").Value = 47681.695
$ws.Range("K134").Value = 143045.085
$ws.Range("M134").Value = -140510.085

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5200.5
$ws.Range("I31").Value = 1109.2559
$ws.Range("J31").Value = 30332.428
$ws.Range("K31").Value = 1109.2559
$ws.Range("L31").Value = 30332.428
$ws.Range("M31").Value = -814.2559000000001
$ws.Range("N31").Value = -30922.428

$ws.Range("H34").Value = 5200.5
$ws.Range("I34").Value = 1109.2559
$ws.Range("J34").Value = 30332.428
$ws.Range("K34").Value = 1109.2559
$ws.Range("L34").Value = 30332.428
$ws.Range("M34").Value = -907.2559000000001
$ws.Range("N34").Value = -30736.428

$ws.Range("H99").Value = 7562780.5
$ws.Range("I99").Value = 5089272
$ws.Range("J99").Value = 13334300
$ws.Range("K99").Value = 5089272
$ws.Range("L99").Value = 13334300
$ws.Range("M99").Value = -5087774
$ws.Range("N99").Value = -13337296

$ws.Range("H126").Value = 7562780.5
$ws.Range("I126").Value = 5089272
$ws.Range("J126").Value = 13334300
$ws.Range("K126").Value = 15267816
$ws.Range("L126").Value = 40002900
$ws.Range("M126").Value = -15265346
$ws.Range("N126").Value = -40007840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1842.5714
$ws.Range("I57").Value = 1842.5714
$ws.Range("K57").Value = 5527.7142
$ws.Range("M57").Value = -4968.7142

$ws.Range("H68").Value = 2116.3125
$ws.Range("I68").Value = 4023.5
$ws.Range("K68").Value = 12070.5
$ws.Range("M68").Value = -11259.5

$ws.Range("H71").Value = 2116.3125
$ws.Range("I71").Value = 4023.5
$ws.Range("K71").Value = 36211.5
$ws.Range("M71").Value = -32155.5

$ws.Range("H137").Value = 4923.5386
$ws.Range("I137").Value = 4600.2
$ws.Range("J137").Value = 5125.625
$ws.Range("K137").Value = 13800.6
$ws.Range("L137").Value = 15376.875
$ws.Range("M137").Value = -8700.599999999999
$ws.Range("N137").Value = -25576.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7124.5
$ws.Range("J70").Value = 9549.799999999999
$ws.Range("L70").Value = 9549.799999999999
$ws.Range("N70").Value = -10089.8

$ws.Range("H73").Value = 7124.5
$ws.Range("J73").Value = 9549.799999999999
$ws.Range("L73").Value = 9549.799999999999
$ws.Range("N73").Value = -11421.8

$ws.Range("H109").Value = 59164.168
$ws.Range("J109").Value = 59164.168
$ws.Range("L109").Value = 59164.168
$ws.Range("N109").Value = -61244.168

$ws.Range("H122").Value = 4332074
$ws.Range("I122").Value = 4950684.5
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 14852053.5
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -14849603.5
$ws.Range("N122").Value = -10300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6537547
$ws.Range("I40").Value = 1399
$ws.Range("J40").Value = 19609842
$ws.Range("K40").Value = 1399
$ws.Range("L40").Value = 19609842
$ws.Range("M40").Value = -1263
$ws.Range("N40").Value = -19610114

$ws.Range("H55").Value = 1382.1428
$ws.Range("J55").Value = 1802.5834
$ws.Range("L55").Value = 1802.5834
$ws.Range("N55").Value = -2148.5834

$ws.Range("H61").Value = 2075
$ws.Range("I61").Value = 2087.5833
$ws.Range("K61").Value = 2087.5833
$ws.Range("M61").Value = -1885.5833

$ws.Range("H93").Value = 83342690
$ws.Range("I93").Value = 111123170
$ws.Range("J93").Value = 1249.6666
$ws.Range("K93").Value = 111123170
$ws.Range("L93").Value = 1249.6666
$ws.Range("M93").Value = -111121922
$ws.Range("N93").Value = -3745.6666

$ws.Range("H100").Value = 2569.7
$ws.Range("I100").Value = 2149.9375
$ws.Range("J100").Value = 4248.75
$ws.Range("K100").Value = 2149.9375
$ws.Range("L100").Value = 4248.75
$ws.Range("M100").Value = -1608.9375
$ws.Range("N100").Value = -5330.75

$ws.Range("H113").Value = 2075
$ws.Range("I113").Value = 2087.5833
$ws.Range("K113").Value = 2087.5833
$ws.Range("M113").Value = 82.41670000000022

$ws.Range("H132").Value = 3679111.5
$ws.Range("I132").Value = 5889
$ws.Range("J132").Value = 6985011.5
$ws.Range("K132").Value = 17667
$ws.Range("L132").Value = 20955034.5
$ws.Range("M132").Value = -15137
$ws.Range("N132").Value = -20960094.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 85000
$ws.Range("J60").Value = 85000
$ws.Range("L60").Value = 85000
$ws.Range("N60").Value = -86644

$ws.Range("H96").Value = 1198.8
$ws.Range("I96").Value = 998.5
$ws.Range("J96").Value = 1332.3334
$ws.Range("K96").Value = 998.5
$ws.Range("L96").Value = 1332.3334
$ws.Range("M96").Value = 374.5
$ws.Range("N96").Value = -4078.3334
